$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for row 6 (analysis #05 / cpx #5 / core / eruption III)
$ws.Range("A6").Value = "analysis #05"
$ws.Range("B6").Value = "cpx #5"
$ws.Range("C6").Value = "core"
$ws.Range("D6").Value = "eruption III"
$ws.Range("E6").Value = 53.1158
$ws.Range("F6").Value = 0.3346
$ws.Range("G6").Value = 1.9724
$ws.Range("H6").Value = 0.5675
$ws.Range("I6").Value = 3.157
$ws.Range("J6").Value = 0.0791
$ws.Range("L6").Value = 17.2283
$ws.Range("M6").Value = 23.9541
$ws.Range("N6").Value = 0.1405

# New data for row 7 (analysis #06 / cpx #6 / rim / eruption III)
$ws.Range("A7").Value = "analysis #06"
$ws.Range("B7").Value = "cpx #6"
$ws.Range("C7").Value = "rim"
$ws.Range("D7").Value = "eruption III"
$ws.Range("E7").Value = 52.6375
$ws.Range("F7").Value = 0.3455
$ws.Range("G7").Value = 2.3905
$ws.Range("H7").Value = 0.6426
$ws.Range("I7").Value = 3.2407
$ws.Range("J7").Value = 0.0874
$ws.Range("L7").Value = 16.5586
$ws.Range("M7").Value = 23.9066
$ws.Range("N7").Value = 0.0686

# New data for row 8 (analysis #07 / cpx #7 / core / eruption IV)
$ws.Range("A8").Value = "analysis #07"
$ws.Range("B8").Value = "cpx #7"
$ws.Range("C8").Value = "core"
$ws.Range("D8").Value = "eruption IV"
$ws.Range("E8").Value = 49.81
$ws.Range("F8").Value = 0.42
$ws.Range("G8").Value = 2.29
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 8.3
$ws.Range("J8").Value = 0.61
$ws.Range("L8").Value = 12.93
$ws.Range("M8").Value = 22.11
$ws.Range("N8").Value = 0.59
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 97.4

# New data for row 9 (analysis #08 / cpx #8 / rim / eruption IV)
$ws.Range("A9").Value = "analysis #08"
$ws.Range("B9").Value = "cpx #8"
$ws.Range("C9").Value = "rim"
$ws.Range("D9").Value = "eruption IV"
$ws.Range("E9").Value = 51.98
$ws.Range("F9").Value = 0.41
$ws.Range("G9").Value = 1.93
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 2.78
$ws.Range("J9").Value = 0.07
$ws.Range("L9").Value = 17
$ws.Range("M9").Value = 23.78
$ws.Range("N9").Value = 0.2
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 98.92

# Rename duplicate header label so it becomes distinct ("notes " -> "notes.1")
$ws.Range("D1").Value = "notes.1"

# Update selection to C1 (matches the diff's sheetView selection change)
$ws.Range("C1").Select()

# Disable iterative calculation (calcPr iterate attribute removed)
$excel.Iteration = $false
